# Update scripts with new TPM (transcripts-per-million) values.
#
# The ligand/receptor edge table (NATMI output) was regenerated from an
# updated expression matrix. All numeric statistics in data rows 2-7 are
# refreshed, and two brand-new target-cluster rows are appended for a
# previously-missing cluster ("Resolving-Mac"), giving the sheet nine
# data rows (rows 2-9) in total instead of six.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Wnt1/Fzd4 -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.0960827240265261
$ws.Range("J2").Value = 0.09608272402652611
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.524618
$ws.Range("N2").Value = 58.573854
$ws.Range("O2").Value = 0.4154885426712971
$ws.Range("P2").Value = 0.4539723485554654
$ws.Range("Q2").Value = 0.529566214014
$ws.Range("R2").Value = 4.766095926126
$ws.Range("S2").Value = 0.03992127098166975
$ws.Range("T2").Value = 0.0436188998819287

# Row 3: ECs -> Wnt1/Fzd4 -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.027123
$ws.Range("H3").Value = 0.081369
$ws.Range("I3").Value = 0.0960827240265261
$ws.Range("J3").Value = 0.09608272402652611
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.24435933333334
$ws.Range("N3").Value = 45.73307800000001
$ws.Range("O3").Value = 0.324403614112412
$ws.Range("P3").Value = 0.3544508583357054
$ws.Range("Q3").Value = 0.413472758198
$ws.Range("R3").Value = 3.721254823782
$ws.Range("S3").Value = 0.03116958292797055
$ws.Range("T3").Value = 0.03405660400243488

# Row 4: ECs -> Wnt1/Fzd4 -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.027123
$ws.Range("H4").Value = 0.081369
$ws.Range("I4").Value = 0.0960827240265261
$ws.Range("J4").Value = 0.09608272402652611
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.9507005
$ws.Range("N4").Value = 23.901401
$ws.Range("O4").Value = 0.2543137660693869
$ws.Range("P4").Value = 0.1852460510065796
$ws.Range("Q4").Value = 0.3241388496615
$ws.Range("R4").Value = 1.944833097969
$ws.Range("S4").Value = 0.02443515940139142
$ws.Range("T4").Value = 0.01779894519586897

# Row 5: ECs -> Wnt1/Fzd4 -> Resolving-Mac (new)
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.027123
$ws.Range("H5").Value = 0.081369
$ws.Range("I5").Value = 0.0960827240265261
$ws.Range("J5").Value = 0.09608272402652611
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.272275
$ws.Range("N5").Value = 0.816825
$ws.Range("O5").Value = 0.005794077146903843
$ws.Range("P5").Value = 0.006330742102249548
$ws.Range("Q5").Value = 0.007384914824999999
$ws.Range("R5").Value = 0.066464233425
$ws.Range("S5").Value = 0.0005567107154943637
$ws.Range("T5").Value = 0.0006082749462935531

# Row 6: FAPs -> Wnt1/Fzd4 -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt1"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.255165
$ws.Range("H6").Value = 0.7654949999999999
$ws.Range("I6").Value = 0.9039172759734738
$ws.Range("J6").Value = 0.9039172759734738
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.524618
$ws.Range("N6").Value = 58.573854
$ws.Range("O6").Value = 0.4154885426712971
$ws.Range("P6").Value = 0.4539723485554654
$ws.Range("Q6").Value = 4.981999151969999
$ws.Range("R6").Value = 44.83799236773
$ws.Range("S6").Value = 0.3755672716896273
$ws.Range("T6").Value = 0.4103534486735367

# Row 7: FAPs -> Wnt1/Fzd4 -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt1"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.255165
$ws.Range("H7").Value = 0.7654949999999999
$ws.Range("I7").Value = 0.9039172759734738
$ws.Range("J7").Value = 0.9039172759734738
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 15.24435933333334
$ws.Range("N7").Value = 45.73307800000001
$ws.Range("O7").Value = 0.324403614112412
$ws.Range("P7").Value = 0.3544508583357054
$ws.Range("Q7").Value = 3.88982694929
$ws.Range("R7").Value = 35.00844254361
$ws.Range("S7").Value = 0.2932340311844415
$ws.Range("T7").Value = 0.3203942543332705

# Row 8: FAPs -> Wnt1/Fzd4 -> MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt1"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.255165
$ws.Range("H8").Value = 0.7654949999999999
$ws.Range("I8").Value = 0.9039172759734738
$ws.Range("J8").Value = 0.9039172759734738
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 11.9507005
$ws.Range("N8").Value = 23.901401
$ws.Range("O8").Value = 0.2543137660693869
$ws.Range("P8").Value = 0.1852460510065796
$ws.Range("Q8").Value = 3.0494004930825
$ws.Range("R8").Value = 18.296402958495
$ws.Range("S8").Value = 0.2298786066679955
$ws.Range("T8").Value = 0.1674471058107107

# Row 9: FAPs -> Wnt1/Fzd4 -> Resolving-Mac (new)
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt1"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.255165
$ws.Range("H9").Value = 0.7654949999999999
$ws.Range("I9").Value = 0.9039172759734738
$ws.Range("J9").Value = 0.9039172759734738
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.272275
$ws.Range("N9").Value = 0.816825
$ws.Range("O9").Value = 0.005794077146903843
$ws.Range("P9").Value = 0.006330742102249548
$ws.Range("Q9").Value = 0.06947505037499999
$ws.Range("R9").Value = 0.625275453375
$ws.Range("S9").Value = 0.005237366431409479
$ws.Range("T9").Value = 0.005722467155955994

